$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column at T ("Approx Color"), shifting the existing
#    Constellation / HD / Bayer / Flamsteed / Hipparcos / Gaia / HR / WISE
#    columns one place to the right (T->U, U->V, ... AA->AB).
# ---------------------------------------------------------------------------
$ws.Columns("T:T").Insert()

# ---------------------------------------------------------------------------
# 2. Add star #13 - Ross 248 (row 14)
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = 13
$ws.Range("D14").Value = 3.1549282999999999
$ws.Range("E14").Value = 0.16
$ws.Range("F14").Value = 0.13600000000000001
$ws.Range("G14").Value = 23.698378999999999
$ws.Range("H14").Value = 44.174923999999997
$ws.Range("N14").Value = 12.3
$ws.Range("O14").Value = 14.79
$ws.Range("P14").Value = [double]"1.8E-3"
$ws.Range("R14").Value = 1.92

$ws.Range("C14").Value = "Gl 905"
$ws.Range("U14").Value = "Andromeda"
$ws.Range("Q14").Value = "M6V"
$ws.Range("B14").Value = "Ross 248"
$ws.Range("W14").Value = "HH Andromedae"

# ---------------------------------------------------------------------------
# 3. Add star #14 - Epsilon Eridani (row 15)
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = 14
$ws.Range("D15").Value = 3.2116495999999999
$ws.Range("E15").Value = 0.73499999999999999
$ws.Range("F15").Value = 0.82
$ws.Range("G15").Value = 3.548848
$ws.Range("H15").Value = -9.4582619999999995
$ws.Range("N15").Value = 3.7360000000000002
$ws.Range("O15").Value = 6.19
$ws.Range("P15").Value = 0.34
$ws.Range("R15").Value = 0.88700000000000001
$ws.Range("V15").Value = 22049
$ws.Range("X15").Value = 18
$ws.Range("Y15").Value = 16537
$ws.Range("AA15").Value = 1084

$ws.Range("W15").Value = "Epsilon Eridani"
$ws.Range("B15").Value = "Epsilon Eridani;Ran"
$ws.Range("U15").Value = "Eridanus"
$ws.Range("Q15").Value = "K2V"
$ws.Range("C15").Value = "Gl 144"

# ---------------------------------------------------------------------------
# 4. Header for the newly inserted column
# ---------------------------------------------------------------------------
$ws.Range("T1").Value = "Approx Color"

# ---------------------------------------------------------------------------
# 5. Fill in the "Approx Color" column for every star row, grouped by color
# ---------------------------------------------------------------------------
foreach ($c in @("T3","T6","T7","T8","T13","T14","T16")) {
    $ws.Range($c).Value = "red"
}
foreach ($c in @("T5","T15")) {
    $ws.Range($c).Value = "orange"
}
foreach ($c in @("T11")) {
    $ws.Range($c).Value = "blue"
}
foreach ($c in @("T2","T4","T9","T10","T12")) {
    $ws.Range($c).Value = "white"
}

# ---------------------------------------------------------------------------
# 6. Add star #15 - Lacaille 9352 (row 16)
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = 15
$ws.Range("D16").Value = 3.2745028999999999
$ws.Range("E16").Value = 0.503
$ws.Range("F16").Value = 0.45900000000000002
$ws.Range("G16").Value = 23.097531
$ws.Range("H16").Value = -35.853073000000002
$ws.Range("N16").Value = 7.34
$ws.Range("O16").Value = 9.8000000000000007
$ws.Range("P16").Value = 0.33
$ws.Range("R16").Value = 1.5
$ws.Range("V16").Value = 217987
$ws.Range("Y16").Value = 114046

$ws.Range("B16").Value = "Lacaille 9352"
$ws.Range("C16").Value = "Gl 887"
$ws.Range("Q16").Value = "M0.5V"
$ws.Range("U16").Value = "Piscus Austrinus"

# ---------------------------------------------------------------------------
# 7. Update the current selection (so the workbook reopens focused on D25,
#    matching the author's last cursor position) and stop the Enter/Tab
#    key presses used while typing the new rows from bubbling up further.
# ---------------------------------------------------------------------------
[void]$ws.Range("D25").Select()
